$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-26 04:46:44"

$wsZhCn.Range("H2").Value = "2016-08-26 04:46:40"
$wsZhCn.Range("K2").Value = "2016-08-26 04:46:56"

$wsDeDe.Range("H2").Value = "2016-08-26 04:46:44"
$wsDeDe.Range("K2").Value = "2016-08-26 04:47:07"
